$wb = $excel.ActiveWorkbook

# Hunk -1477,25 +1477,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4020519.8
$ws.Range("J17").Value = 6699999.5
$ws.Range("L17").Value = 20099998.5
$ws.Range("N17").Value = -20100334.5

# Hunk -7232,22 +7232,22 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3226.1538
$ws.Range("I132").Value = 3228.0222
$ws.Range("K132").Value = 9684.0666
$ws.Range("M132").Value = -7154.0666

# Hunk -7379,25 +7379,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2003.3077
$ws.Range("I135").Value = 1464.5
$ws.Range("J135").Value = 3799.3333
$ws.Range("K135").Value = 13180.5
$ws.Range("L135").Value = 34193.9997
$ws.Range("M135").Value = -10645.5
$ws.Range("N135").Value = -39263.9997

# Hunk -7532,25 +7532,25 @@ (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3017.0715
$ws.Range("I138").Value = 2167.85
$ws.Range("J138").Value = 3789.0908
$ws.Range("K138").Value = 6503.549999999999
$ws.Range("L138").Value = 11367.2724
$ws.Range("M138").Value = -1363.549999999999
$ws.Range("N138").Value = -21647.2724

# Hunk -12024,25 +12024,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3469.5
$ws.Range("J88").Value = 4850
$ws.Range("L88").Value = 4850
$ws.Range("N88").Value = -5662

# Hunk -12168,25 +12168,25 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3469.5
$ws.Range("J91").Value = 4850
$ws.Range("L91").Value = 4850
$ws.Range("N91").Value = -7658

# Hunk -13993,22 +13993,22 @@ (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 102490
$ws.Range("J128").Value = 102490
$ws.Range("L128").Value = 102490
$ws.Range("N128").Value = -112450

# Hunk -23163,25 +23163,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3850.25
$ws.Range("I31").Value = 2083.25
$ws.Range("J31").Value = 5617.25
$ws.Range("K31").Value = 2083.25
$ws.Range("L31").Value = 5617.25
$ws.Range("M31").Value = -1788.25
$ws.Range("N31").Value = -6207.25

# Hunk -23316,25 +23316,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3850.25
$ws.Range("I34").Value = 2083.25
$ws.Range("J34").Value = 5617.25
$ws.Range("K34").Value = 2083.25
$ws.Range("L34").Value = 5617.25
$ws.Range("M34").Value = -1881.25
$ws.Range("N34").Value = -6021.25

# Hunk -25313,22 +25313,22 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 104640.664
$ws.Range("J75").Value = 104640.664
$ws.Range("L75").Value = 104640.664
$ws.Range("N75").Value = -106636.664

# Hunk -25460,22 +25460,22 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H78").Value = 104640.664
$ws.Range("J78").Value = 104640.664
$ws.Range("L78").Value = 313921.992
$ws.Range("N78").Value = -323905.992

# Hunk -26547,22 +26547,22 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 110750
$ws.Range("J100").Value = 110750
$ws.Range("L100").Value = 110750
$ws.Range("N100").Value = -112914

# Hunk -28136,22 +28136,22 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1932.2354
$ws.Range("I132").Value = 1523.2
$ws.Range("K132").Value = 4569.6
$ws.Range("M132").Value = -2039.6

# Hunk -28234,25 +28234,25 @@ (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2518.2334
$ws.Range("I134").Value = 1318.0454
$ws.Range("J134").Value = 5818.75
$ws.Range("K134").Value = 3954.1362
$ws.Range("L134").Value = 17456.25
$ws.Range("M134").Value = -1419.1362
$ws.Range("N134").Value = -22526.25

# Hunk -34237,25 +34237,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1184.6522
$ws.Range("I113").Value = 579.75
$ws.Range("J113").Value = 1312
$ws.Range("K113").Value = 1739.25
$ws.Range("L113").Value = 3936
$ws.Range("M113").Value = 430.75
$ws.Range("N113").Value = -8276

# Hunk -34580,25 +34580,25 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 24214.445
$ws.Range("I120").Value = 22982.5
$ws.Range("J120").Value = 25200
$ws.Range("K120").Value = 68947.5
$ws.Range("L120").Value = 75600
$ws.Range("M120").Value = -64109.5
$ws.Range("N120").Value = -85276

# Hunk -35134,22 +35134,22 @@ (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1707.9056
$ws.Range("I131").Value = 1216.5385
$ws.Range("K131").Value = 3649.6155
$ws.Range("M131").Value = 1390.3845

# Hunk -41932,22 +41932,22 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 134729
$ws.Range("J128").Value = 134729
$ws.Range("L128").Value = 134729
$ws.Range("N128").Value = -144689

# Hunk -42125,22 +42125,22 @@ (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2582.125
$ws.Range("I132").Value = 2420.9333
$ws.Range("K132").Value = 7262.7999
$ws.Range("M132").Value = -4732.7999

# Hunk -43303,22 +43303,19 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Hunk -44512,22 +44509,19 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

# Hunk -44904,25 +44898,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3112.4546
$ws.Range("J46").Value = 3370.889
$ws.Range("L46").Value = 3370.889
$ws.Range("N46").Value = -3746.889

# Hunk -45991,25 +45985,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3074.4285
$ws.Range("I68").Value = 1261
$ws.Range("J68").Value = 3799.8
$ws.Range("K68").Value = 1261
$ws.Range("L68").Value = 3799.8
$ws.Range("M68").Value = -512
$ws.Range("N68").Value = -5297.8

# Hunk -46135,25 +46129,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3074.4285
$ws.Range("I71").Value = 1261
$ws.Range("J71").Value = 3799.8
$ws.Range("K71").Value = 6305
$ws.Range("L71").Value = 18999
$ws.Range("M71").Value = -2561
$ws.Range("N71").Value = -26487

# Hunk -46665,25 +46659,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2921.3914
$ws.Range("J82").Value = 3331.6155
$ws.Range("L82").Value = 3331.6155
$ws.Range("N82").Value = -4053.6155

# Hunk -46809,25 +46803,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2921.3914
$ws.Range("J85").Value = 3331.6155
$ws.Range("L85").Value = 3331.6155
$ws.Range("N85").Value = -5827.6155

# Hunk -48628,25 +48622,25 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2904.0833
$ws.Range("I122").Value = 3049.8096
$ws.Range("J122").Value = 1884
$ws.Range("K122").Value = 9149.4288
$ws.Range("L122").Value = 5652
$ws.Range("M122").Value = -6699.4288
$ws.Range("N122").Value = -10552

# Hunk -49271,22 +49265,22 @@ (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 49979
$ws.Range("J135").Value = 49979
$ws.Range("L135").Value = 49979
$ws.Range("N135").Value = -60119

# Hunk -53579,25 +53573,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2309.1667
$ws.Range("I81").Value = 1397.5
$ws.Range("J81").Value = 5500
$ws.Range("K81").Value = 2795
$ws.Range("L81").Value = 11000
$ws.Range("M81").Value = -1734
$ws.Range("N81").Value = -13122

# Hunk -53726,25 +53720,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2309.1667
$ws.Range("I84").Value = 1397.5
$ws.Range("J84").Value = 5500
$ws.Range("K84").Value = 13975
$ws.Range("L84").Value = 55000
$ws.Range("M84").Value = -8671
$ws.Range("N84").Value = -65608

# Hunk -55597,25 +55591,25 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3231.889
$ws.Range("I122").Value = 2512.2666
$ws.Range("J122").Value = 6830
$ws.Range("K122").Value = 7536.7998
$ws.Range("L122").Value = 20490
$ws.Range("M122").Value = -5086.7998
$ws.Range("N122").Value = -25390

# Hunk -56280,22 +56274,22 @@ (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21506.303
$ws.Range("I136").Value = 1614.5454
$ws.Range("K136").Value = 4843.6362
$ws.Range("M136").Value = -2293.6362
